$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert 3 new rows before the old row 40 ("Total" row), pushing it down to row 43 ---
$ws.Range("A40:A42").EntireRow.Insert()

# --- 2. Build the new "Dimensions" rows (40,41,42) by copying format+content from the
#        existing "Projected Area" rows (38,39) which have the matching per-row styling,
#        then overwrite the text/values with the new Dimensions data ---
$ws.Range("A38:E38").Copy($ws.Range("A40"))
$ws.Range("A39:E39").Copy($ws.Range("A41"))
$ws.Range("A41:E41").Copy($ws.Range("A42"))

# Row 40: Dimensions / X / DIM_X
$ws.Range("A40").Value = "Dimensions"
$ws.Range("B40").Value = "X"
$ws.Range("C40").Value = "DIM_X"
$ws.Range("D40").Value = 806.45
$ws.Range("E40").Value = "mm"

# Row 41: Y / DIM_Y
$ws.Range("B41").Value = "Y"
$ws.Range("C41").Value = "DIM_Y"
$ws.Range("D41").Value = 419.1
$ws.Range("E41").Value = "mm"

# Row 42: Z / DIM_Z
$ws.Range("B42").Value = "Z"
$ws.Range("C42").Value = "DIM_Z"
$ws.Range("D42").Value = 330.2
$ws.Range("E42").Value = "mm"

# --- 3. Add new column F (second set of readings) for the Projected Area rows (37-39)
#        and the new Dimensions rows (40-42), copying the number format from column D ---
$ws.Range("D37").Copy($ws.Range("F37"))
$ws.Range("D38").Copy($ws.Range("F38"))
$ws.Range("D39").Copy($ws.Range("F39"))
$ws.Range("D40").Copy($ws.Range("F40"))
$ws.Range("D41").Copy($ws.Range("F41"))
$ws.Range("D42").Copy($ws.Range("F42"))

# --- 4. Update the values (Projected Area D column corrected + new F column) ---
$ws.Range("D37").Value = 0.109066
$ws.Range("F37").Value = 0.230487

$ws.Range("D38").Value = 0.224288
$ws.Range("F38").Value = 0.224288

$ws.Range("D39").Value = 0.230487
$ws.Range("F39").Value = 0.109066

$ws.Range("F40").Value = 330.2
$ws.Range("F41").Value = 419.1
$ws.Range("F42").Value = 806.45

# --- 5. Update sheet view to match the saved selection/scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("F33").Select()
